# Fill in the previously-missing F1/Precision/Recall/Accuracy result strings
# for the "Count Vectorizer + TFIDF + ngram(2) + POS" block (rows 51-55) on
# the "Lucene" sheet. These values were scraped/added later for the
# Logistic Regression, Multinomial Naive Bayes, Support Vector Machines,
# Decision Tree and Random Forest rows, mirroring the other blocks already
# present in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lucene")

# Row 51 - Logistic Regression
$ws.Range("C51").Value = "0.381 0.724 0.637 0.437 0.722 "
$ws.Range("D51").Value = "0.298 0.646 0.606 0.373 0.647"
$ws.Range("E51").Value = "0.239 0.649 0.490 0.285 0.575 "
$ws.Range("F51").Value = "0.860 0.761 0.816 0.861 0.937"

# Row 52 - Multinomial Naive Bayes
$ws.Range("C52").Value = "0.420 0.750 0.670 0.519 0.746"
$ws.Range("D52").Value = "0.370 0.585 0.534 0.430 0.702"
$ws.Range("E52").Value = "0.270 0.792 0.547 0.359 0.604"
$ws.Range("F52").Value = "0.874 0.739 0.794 0.870 0.945 "

# Row 53 - Support Vector Machines
$ws.Range("C53").Value = "0.448 0.750 0.643 0.559 0.709"
$ws.Range("D53").Value = "0.406 0.657 0.560 0.456 0.732 "
$ws.Range("E53").Value = "0.293 0.694 0.505 0.398 0.556"
$ws.Range("F53").Value = "0.880 0.774 0.802 0.875 0.945 "

# Row 54 - Decision Tree
$ws.Range("C54").Value = "0.287 0.671 0.385 0.229 0.641"
$ws.Range("D54").Value = "0.204 0.539 0.491 0.287 0.852 "
$ws.Range("E54").Value = "0.170 0.623 0.243 0.130 0.473"
$ws.Range("F54").Value = "0.842 0.692 0.776 0.861 0.948 "

# Row 55 - Random Forest
$ws.Range("C55").Value = "0.243 0.762 0.571 0.367 0.707"
$ws.Range("D55").Value = "0.456 0.617 0.718 0.660 0.781"
$ws.Range("E55").Value = "0.139 0.769 0.407 0.225 0.551"
$ws.Range("F55").Value = "0.891 0.760 0.833 0.897 0.949 "

# The author had last clicked on C56 before saving, update the sheet's
# active selection accordingly.
$ws.Range("C56").Select()
